# Apply the cell-value edits described by the diff.
# Only the text of specific B:G cells changes; column A values are untouched.
# The numeric-looking values (e.g. "3", "14", "36") must remain stored as TEXT
# (they were t="inlineStr" before and stay that type after), so we force a
# text number format before writing them, then restore the default style so
# no extra style index gets attached to the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: append a trailing "d" to B2:E2 (F2 and G2 stay the same, already text)
$ws.Range("B2").Value = "егорdsa2d"
$ws.Range("C2").Value = "губин2d"
$ws.Range("D2").Value = "выфывфы2d"
$ws.Range("E2").Value = "выфв2d"

# Rows 3-7, columns B:G: new values that look numeric, so force text storage.
$numericRange = $ws.Range("B3:G7")
$numericRange.NumberFormat = "@"

# Row 3: all of B3:G3 become "3"
$ws.Range("B3:G3").Value = "3"

# Row 4: B4:E4 and G4 become "14", F4 becomes "41"
$ws.Range("B4").Value = "14"
$ws.Range("C4").Value = "14"
$ws.Range("D4").Value = "14"
$ws.Range("E4").Value = "14"
$ws.Range("F4").Value = "41"
$ws.Range("G4").Value = "14"

# Row 5: all of B5:G5 become "1"
$ws.Range("B5:G5").Value = "1"

# Row 6: all of B6:G6 become "36"
$ws.Range("B6:G6").Value = "36"

# Row 7: all of B7:G7 become "3"
$ws.Range("B7:G7").Value = "3"

# Restore default styling so no leftover per-cell style index remains.
$numericRange.Style = "Normal"
